$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the formula in B9: from (B3*B6 - B3*B7) to (B3-B8)
$ws.Range("B9").Formula = "=B3-B8"

# Update the sheet view: zoom to 130% and change selection to C9
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("C9").Select() | Out-Null
